$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (A4/B4) had been copy-pasted from row 3 (user_id_2 / 用户2).
# Fix it so it refers to the 3rd user, matching C4 (user_id_3@1.c).
$ws.Range("A4").Value = "user_id_3"
$ws.Range("B4").Value = "用户3"

$ws.Range("C17").Select()
